$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 25.63000000000057
$ws.Range("H2").Value = [double]"1.033913669856723e-05"
$ws.Range("I2").Value = [double]"1.033913669856723e-05"
$ws.Range("L2").Value = 47.44987279384012
$ws.Range("M2").Value = "[26.018595720753268, 68.88114986692698]"
$ws.Range("N2").Value = [double]"5.429484279217967e-05"
$ws.Range("O2").Value = [double]"5.429484279217967e-05"
$ws.Range("P2").Value = 1.905710858934349
$ws.Range("Q2").Value = "[1.3522370781217345, 2.4591846397469643]"
$ws.Range("R2").Value = [double]"1.272621563686016e-08"
$ws.Range("S2").Value = [double]"1.272621563686016e-08"
$ws.Range("T2").Value = 64.4500992842018
$ws.Range("U2").Value = "[51.65299297278321, 77.24720559562039]"
$ws.Range("V2").Value = [double]"3.319566843629218e-13"
$ws.Range("W2").Value = [double]"3.319566843629218e-13"
$ws.Range("X2").Value = 17.85633633633673
$ws.Range("Y2").Value = 15.59863863863899
$ws.Range("Z2").Value = 20.11403403403447

# Row 3
$ws.Range("F3").Value = 25.63000000000057
$ws.Range("H3").Value = 0.0078217854411754
$ws.Range("I3").Value = 0.0078217854411754
$ws.Range("L3").Value = 32.42928658882848
$ws.Range("M3").Value = "[7.478654160769182, 57.37991901688777]"
$ws.Range("N3").Value = 0.01201233096790477
$ws.Range("O3").Value = 0.01201233096790477
$ws.Range("P3").Value = 1.33965812855781
$ws.Range("Q3").Value = "[0.32076321388003937, 2.35855304323558]"
$ws.Range("R3").Value = 0.0111195262563677
$ws.Range("S3").Value = 0.0111195262563677
$ws.Range("T3").Value = 56.48936424907076
$ws.Range("U3").Value = "[42.14969906206289, 70.82902943607863]"
$ws.Range("V3").Value = [double]"4.296523137270469e-10"
$ws.Range("W3").Value = [double]"4.296523137270469e-10"
$ws.Range("X3").Value = 20.16534534534579
$ws.Range("Y3").Value = 16.00912912912948
$ws.Range("Z3").Value = 24.3215615615621

# Row 4
$ws.Range("F4").Value = 25.63000000000057
$ws.Range("H4").Value = [double]"4.095478231713745e-05"
$ws.Range("I4").Value = [double]"4.095478231713745e-05"
$ws.Range("L4").Value = 49.98679813858119
$ws.Range("M4").Value = "[22.78658439485318, 77.18701188230919]"
$ws.Range("N4").Value = 0.000582530909898793
$ws.Range("O4").Value = 0.000582530909898793
$ws.Range("P4").Value = 0.9371317425122703
$ws.Range("Q4").Value = "[0.3962369112635775, 1.478026573760963]"
$ws.Range("R4").Value = 0.00109459314563054
$ws.Range("S4").Value = 0.00109459314563054
$ws.Range("T4").Value = 63.43523018156016
$ws.Range("U4").Value = "[49.42491651201276, 77.44554385110756]"
$ws.Range("V4").Value = [double]"8.602230039400638e-12"
$ws.Range("W4").Value = [double]"8.602230039400638e-12"
$ws.Range("X4").Value = 21.80730730730779
$ws.Range("Y4").Value = 19.60092092092136
$ws.Range("Z4").Value = 24.01369369369423

# Row 5
$ws.Range("B5").Value = 1
$ws.Range("F5").Value = 25.63000000000057
$ws.Range("H5").Value = 0.005182374329031791
$ws.Range("I5").Value = 0.005182374329031791
$ws.Range("L5").Value = 29.83292929366689
$ws.Range("M5").Value = "[6.669967382484501, 52.99589120484927]"
$ws.Range("N5").Value = 0.01275508150479343
$ws.Range("O5").Value = 0.01275508150479343
$ws.Range("P5").Value = 0.8239211964369622
$ws.Range("Q5").Value = "[0.14465791998511435, 1.50318447288881]"
$ws.Range("R5").Value = 0.01855281250539043
$ws.Range("S5").Value = 0.01855281250539043
$ws.Range("T5").Value = 54.3233483613021
$ws.Range("U5").Value = "[42.14553037890926, 66.50116634369493]"
$ws.Range("V5").Value = [double]"1.332178811708218e-11"
$ws.Range("W5").Value = [double]"1.332178811708218e-11"
$ws.Range("X5").Value = 22.2691091091096
$ws.Range("Y5").Value = 19.49829829829873
$ws.Range("Z5").Value = 25.03991991992048

# Row 6
$ws.Range("B6").Value = 0
$ws.Range("F6").Value = 24.87000000000045
$ws.Range("H6").Value = 0.05240686171093645
$ws.Range("I6").Value = 0.05240686171093645
$ws.Range("L6").Value = 25.19913931821711
$ws.Range("M6").Value = "[-1.368350920264838, 51.766629556699066]"
$ws.Range("N6").Value = 0.06247078234909353
$ws.Range("O6").Value = 0.06247078234909353
$ws.Range("P6").Value = 0.9119738433844251
$ws.Range("Q6").Value = "[-1.1383949355350396, 2.96234262230389]"
$ws.Range("R6").Value = 0.3751036337645686
$ws.Range("S6").Value = 0.3751036337645686
$ws.Range("T6").Value = 50.94846952488105
$ws.Range("U6").Value = "[36.68657537680822, 65.21036367295387]"
$ws.Range("V6").Value = [double]"5.23896237503152e-09"
$ws.Range("W6").Value = [double]"5.23896237503152e-09"
$ws.Range("X6").Value = 21.26024024024062
$ws.Range("Y6").Value = 13.14450450450474
$ws.Range("Z6").Value = 29.3759759759765

# Row 7
$ws.Range("F7").Value = 24.87000000000045
$ws.Range("H7").Value = 0.001890896205892889
$ws.Range("I7").Value = 0.001890896205892889
$ws.Range("L7").Value = 36.49699498692633
$ws.Range("M7").Value = "[11.529586113434583, 61.46440386041808]"
$ws.Range("N7").Value = 0.005107842410848606
$ws.Range("O7").Value = 0.005107842410848606
$ws.Range("P7").Value = 0.3836579616996545
$ws.Range("Q7").Value = "[-0.37107901213573147, 1.1383949355350405]"
$ws.Range("R7").Value = 0.3113853704198866
$ws.Range("S7").Value = 0.3113853704198866
$ws.Range("T7").Value = 63.27134745931208
$ws.Range("U7").Value = "[49.83165325823341, 76.71104166039075]"
$ws.Range("V7").Value = [double]"2.678746113815578e-12"
$ws.Range("W7").Value = [double]"2.678746113815578e-12"
$ws.Range("X7").Value = 23.35141141141183
$ws.Range("Y7").Value = 20.36402402402439
$ws.Range("Z7").Value = 26.33879879879927

# Row 8
$ws.Range("F8").Value = 24.87000000000045
$ws.Range("H8").Value = [double]"2.240865606628617e-06"
$ws.Range("I8").Value = [double]"2.240865606628617e-06"
$ws.Range("L8").Value = 47.49481261621324
$ws.Range("M8").Value = "[27.715948347390906, 67.27367688503557]"
$ws.Range("N8").Value = [double]"1.582572093572487e-05"
$ws.Range("O8").Value = [double]"1.582572093572487e-05"
$ws.Range("P8").Value = 0.3710790121357315
$ws.Range("Q8").Value = "[-0.08176317216550011, 0.823921196436963]"
$ws.Range("R8").Value = 0.105815293138003
$ws.Range("S8").Value = 0.105815293138003
$ws.Range("T8").Value = 60.83380014616964
$ws.Range("U8").Value = "[49.62905658646774, 72.03854370587155]"
$ws.Range("V8").Value = [double]"2.930988785010413e-14"
$ws.Range("W8").Value = [double]"2.930988785010413e-14"
$ws.Range("X8").Value = 23.40120120120162
$ws.Range("Y8").Value = 21.60876876876916
$ws.Range("Z8").Value = 25.19363363363409

# Row 9
$ws.Range("F9").Value = 24.87000000000045
$ws.Range("H9").Value = 0.000293717665778126
$ws.Range("I9").Value = 0.000293717665778126
$ws.Range("L9").Value = 42.43186396442996
$ws.Range("M9").Value = "[19.967685883967405, 64.89604204489251]"
$ws.Range("N9").Value = 0.0004262437963751253
$ws.Range("O9").Value = 0.0004262437963751253
$ws.Range("P9").Value = 0.106921071293347
$ws.Range("Q9").Value = "[-0.5597632555945768, 0.7736053981812709]"
$ws.Range("R9").Value = 0.7481789085074544
$ws.Range("S9").Value = 0.7481789085074544
$ws.Range("T9").Value = 53.86632548145424
$ws.Range("U9").Value = "[40.47882699428543, 67.25382396862305]"
$ws.Range("V9").Value = [double]"2.433873103058204e-10"
$ws.Range("W9").Value = [double]"2.433873103058204e-10"
$ws.Range("X9").Value = 24.44678678678723
$ws.Range("Y9").Value = 21.80792792792832
$ws.Range("Z9").Value = 27.08564564564613
